$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 33   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/26/2026  Through  2/1/2026"

# --- Crime data table updates (rows 15-28) ---

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = '#,##0'
$ws.Range("K15").Value = 0
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L15").Value = 0
$ws.Range("L15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N15").Value = 0
$ws.Range("N15").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = -28.571428571428
$ws.Range("L16").Value = -50
$ws.Range("M16").Value = -66.666666666666
$ws.Range("N16").Value = -93.243243243243

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 8.333333333333
$ws.Range("I17").Value = 18
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = 38.461538461538
$ws.Range("L17").Value = -18.181818181818
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 125

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = -18.181818181818
$ws.Range("L18").Value = -10
$ws.Range("M18").Value = 125
$ws.Range("N18").Value = -75

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 26.666666666666
$ws.Range("I19").Value = 45
$ws.Range("J19").Value = 36
$ws.Range("K19").Value = 25
$ws.Range("L19").Value = 50
$ws.Range("M19").Value = 221.428571428571
$ws.Range("N19").Value = 73.076923076923

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("E20").Value = -100
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G20").Value = 1
$ws.Range("G20").NumberFormat = '#,##0'
$ws.Range("H20").Value = -100
$ws.Range("H20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J20").Value = 1
$ws.Range("J20").NumberFormat = '#,##0'
$ws.Range("K20").Value = -100
$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 21
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -30.434782608695
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = 8.196721311475
$ws.Range("I21").Value = 78
$ws.Range("J21").Value = 70
$ws.Range("K21").Value = 11.428571428571
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 52.941176470588
$ws.Range("N21").Value = -58.064516129032

# Row 22
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = 150

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 45.454545454545
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = 28.571428571428
$ws.Range("L23").Value = 38.461538461538
$ws.Range("M23").Value = 28.571428571428

# Row 24
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 21.739130434782
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = 21.25
$ws.Range("I24").Value = 102
$ws.Range("J24").Value = 88
$ws.Range("K24").Value = 15.909090909090
$ws.Range("L24").Value = 3.030303030303
$ws.Range("M24").Value = 72.881355932203

# Row 25
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 45.454545454545
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = 9.302325581395
$ws.Range("I25").Value = 52
$ws.Range("J25").Value = 45
$ws.Range("K25").Value = 15.555555555555
$ws.Range("L25").Value = -17.460317460317

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -18.181818181818
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -15.151515151515
$ws.Range("I26").Value = 30
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = -11.764705882352
$ws.Range("L26").Value = -26.829268292682
$ws.Range("M26").Value = 11.111111111111

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G27").Value = 1
$ws.Range("G27").NumberFormat = '#,##0'
$ws.Range("H27").Value = 0
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J27").Value = 1
$ws.Range("J27").NumberFormat = '#,##0'
$ws.Range("K27").Value = 0
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L27").Value = -50
$ws.Range("L27").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J28").Value = 4
